$d = $word.ActiveDocument
$d.Content.Find.Execute("505×3=", $true, $false, $false, $false, $false, $true, 1, $false, "101×3=", 2)
$d.Content.Find.Execute("552×8=", $true, $false, $false, $false, $false, $true, 1, $false, "138×3=", 2)
$d.Content.Find.Execute("286×2=", $true, $false, $false, $false, $false, $true, 1, $false, "288×7=", 2)
$d.Content.Find.Execute("757×9=", $true, $false, $false, $false, $false, $true, 1, $false, "170×5=", 2)
$d.Content.Find.Execute("795×3=", $true, $false, $false, $false, $false, $true, 1, $false, "321×9=", 2)
$d.Content.Find.Execute("827×8=", $true, $false, $false, $false, $false, $true, 1, $false, "587×7=", 2)
$d.Content.Find.Execute("255×2=", $true, $false, $false, $false, $false, $true, 1, $false, "169×2=", 2)
$d.Content.Find.Execute("399×8=", $true, $false, $false, $false, $false, $true, 1, $false, "464×8=", 2)
$d.Content.Find.Execute("223×6=", $true, $false, $false, $false, $false, $true, 1, $false, "307×5=", 2)
$d.Content.Find.Execute("621×2=", $true, $false, $false, $false, $false, $true, 1, $false, "757×7=", 2)
$d.Content.Find.Execute("209×4=", $true, $false, $false, $false, $false, $true, 1, $false, "452×2=", 2)
$d.Content.Find.Execute("527×2=", $true, $false, $false, $false, $false, $true, 1, $false, "384×2=", 2)
$d.Content.Find.Execute("192×5=", $true, $false, $false, $false, $false, $true, 1, $false, "131×3=", 2)
$d.Content.Find.Execute("699×8=", $true, $false, $false, $false, $false, $true, 1, $false, "460×2=", 2)
$d.Content.Find.Execute("510×8=", $true, $false, $false, $false, $false, $true, 1, $false, "526×4=", 2)
$d.Content.Find.Execute("764×5=", $true, $false, $false, $false, $false, $true, 1, $false, "545×4=", 2)
$d.Content.Find.Execute("414×3=", $true, $false, $false, $false, $false, $true, 1, $false, "305×7=", 2)
$d.Content.Find.Execute("528×7=", $true, $false, $false, $false, $false, $true, 1, $false, "173×2=", 2)
$d.Content.Find.Execute("156×7=", $true, $false, $false, $false, $false, $true, 1, $false, "134×2=", 2)
$d.Content.Find.Execute("623×4=", $true, $false, $false, $false, $false, $true, 1, $false, "255×9=", 2)
$d.Content.Find.Execute("197×8=", $true, $false, $false, $false, $false, $true, 1, $false, "697×2=", 2)
$d.Content.Find.Execute("367×8=", $true, $false, $false, $false, $false, $true, 1, $false, "650×8=", 2)
$d.Content.Find.Execute("453×9=", $true, $false, $false, $false, $false, $true, 1, $false, "984×4=", 2)
$d.Content.Find.Execute("120×6=", $true, $false, $false, $false, $false, $true, 1, $false, "342×5=", 2)
$d.Content.Find.Execute("257×6=", $true, $false, $false, $false, $false, $true, 1, $false, "378×4=", 2)
